$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

    $ws.Range("E2").Value = 99
    $ws.Range("F2").Value = 72
    $ws.Range("H2").Value = 78
    $ws.Range("F4").Value = 29
    $ws.Range("H4").Value = 41
    $ws.Range("F5").Value = 93
    $ws.Range("H5").Value = 104
    $ws.Range("E6").Value = 45
    $ws.Range("F6").Value = 33
    $ws.Range("H6").Value = 43
    $ws.Range("F7").Value = 18
    $ws.Range("H7").Value = 23
    $ws.Range("E10").Value = 588
    $ws.Range("F10").Value = 298
    $ws.Range("H10").Value = 394
    $ws.Range("E11").Value = 370
    $ws.Range("F11").Value = 203
    $ws.Range("H11").Value = 267
    $ws.Range("E12").Value = 572
    $ws.Range("F12").Value = 319
    $ws.Range("H12").Value = 405
    $ws.Range("E13").Value = 142
    $ws.Range("F13").Value = 78
    $ws.Range("H13").Value = 112
    $ws.Range("F14").Value = 70
    $ws.Range("H14").Value = 104
    $ws.Range("E15").Value = 174
    $ws.Range("E16").Value = 210
    $ws.Range("F16").Value = 106
    $ws.Range("H16").Value = 154
    $ws.Range("F17").Value = 56
    $ws.Range("H17").Value = 80
    $ws.Range("F21").Value = 79
    $ws.Range("H21").Value = 110
    $ws.Range("F22").Value = 95
    $ws.Range("H22").Value = 137
    $ws.Range("E23").Value = 204
    $ws.Range("F23").Value = 100
    $ws.Range("H23").Value = 151
    $ws.Range("F24").Value = 120
    $ws.Range("H24").Value = 150
    $ws.Range("E25").Value = 280
    $ws.Range("F25").Value = 145
    $ws.Range("H25").Value = 205
    $ws.Range("F26").Value = 98
    $ws.Range("H26").Value = 123
    $ws.Range("E27").Value = 338
    $ws.Range("F27").Value = 180
    $ws.Range("H27").Value = 261
    $ws.Range("F28").Value = 83
    $ws.Range("H28").Value = 135
    $ws.Range("E29").Value = 171
    $ws.Range("F29").Value = 101
    $ws.Range("H29").Value = 142
    $ws.Range("F30").Value = 131
    $ws.Range("H30").Value = 183
    $ws.Range("E31").Value = 74
    $ws.Range("F31").Value = 33
    $ws.Range("H31").Value = 61
    $ws.Range("E32").Value = 187
    $ws.Range("F32").Value = 115
    $ws.Range("H32").Value = 153
    $ws.Range("E33").Value = 303
    $ws.Range("F33").Value = 161
    $ws.Range("H33").Value = 250
    $ws.Range("E34").Value = 222
    $ws.Range("F34").Value = 152
    $ws.Range("H34").Value = 191
    $ws.Range("F36").Value = 47
    $ws.Range("H36").Value = 57
    $ws.Range("F37").Value = 88
    $ws.Range("H37").Value = 125
    $ws.Range("F38").Value = 59
    $ws.Range("H38").Value = 75
    $ws.Range("F39").Value = 94
    $ws.Range("H39").Value = 145
    $ws.Range("E40").Value = 268
    $ws.Range("F40").Value = 127
    $ws.Range("H40").Value = 207
    $ws.Range("E41").Value = 397
    $ws.Range("F41").Value = 193
    $ws.Range("H41").Value = 285
    $ws.Range("E42").Value = 390
    $ws.Range("F42").Value = 218
    $ws.Range("H42").Value = 279
    $ws.Range("E43").Value = 124
    $ws.Range("F43").Value = 67
    $ws.Range("H43").Value = 94
    $ws.Range("E44").Value = 317
    $ws.Range("F44").Value = 163
    $ws.Range("H44").Value = 231
    $ws.Range("E45").Value = 150
    $ws.Range("F45").Value = 74
    $ws.Range("H45").Value = 113
    $ws.Range("E46").Value = 334
    $ws.Range("F46").Value = 187
    $ws.Range("H46").Value = 250
    $ws.Range("E47").Value = 465
    $ws.Range("F47").Value = 246
    $ws.Range("H47").Value = 338
    $ws.Range("E48").Value = 220
    $ws.Range("F48").Value = 97
    $ws.Range("H48").Value = 141
    $ws.Range("E49").Value = 290
    $ws.Range("F49").Value = 133
    $ws.Range("H49").Value = 220
    $ws.Range("E50").Value = 247
    $ws.Range("F50").Value = 121
    $ws.Range("H50").Value = 192
    $ws.Range("E51").Value = 243
    $ws.Range("F51").Value = 110
    $ws.Range("H51").Value = 184
